$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = 2
$ws.Range("D4").Value = 0.1
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("D9").Value = 3
$ws.Range("D10").Value = $True
$ws.Range("D13").Value = "[-3000,3000]"
$ws.Range("D14").Value = 2600
$ws.Range("D15").Value = 0
$ws.Range("D18").Value = "[-1000,-1000;1000,1000]"
$ws.Range("D19").Value = "[1600,1600]"
$ws.Range("D20").Value = "[8,0;13,0]"
$ws.Range("D21").Value = 9.81
$ws.Range("C25").Copy($ws.Range("D25")) | Out-Null
$ws.Range("D26").Value = 100
$ws.Range("C28").Copy($ws.Range("D28")) | Out-Null
$ws.Range("D29").Value = 100
$ws.Range("D31").Value = 0.003
$ws.Range("D32").Value = -3
$ws.Range("D33").Value = 10
$ws.Range("D34").Value = 0.5
$ws.Range("D35").Value = 4
$ws.Range("D36").Value = 50
$ws.Range("C38").Copy($ws.Range("D38")) | Out-Null
$ws.Range("D40").Value = 0.5
$ws.Range("D41").Value = "[0.3,0.5]"
$ws.Range("D44").Value = 1000
$ws.Range("D45").Value = 5
$ws.Range("D46").Value = 5
$ws.Range("D47").Value = 15
$ws.Range("D48").Formula = "=-5/12*PI()"
$ws.Range("D49").Formula = "=5/12*PI()"
$ws.Range("D50").Formula = "=11/6*PI()"
$ws.Range("D51").Value = -0.01843
$ws.Range("D52").Value = 0.37819999999999998
$ws.Range("D53").Formula = "=-2.3782"
$ws.Range("D54").Value = 4
$ws.Range("D57").Value = "[-0.5,0.5,-0.5; -0.375,0,0.375]"
$ws.Range("D58").Value = "[-0.5,-0.3,0,0.1,0.2,0.3,0.5,0.3,0.2,0.1,0,-0.3,-0.5;-0.2,-0.1,-0.1,-0.5,-0.5,-0.1,0,0.1,0.5,0.5,0.1,0.1,0.2]"
$ws.Range("D59").Value = "[2, 1.5, 1.5, 0, 0, 1.5, 1.5; 0, 0.5, 0.2, 0.2, -0.2, -0.2, -0.5]"
$ws.Range("D60").Value = "[150;150]"
$ws.Range("C61").Copy($ws.Range("D61")) | Out-Null
$ws.Range("D62").Value = $True
$ws.Range("D63").Value = $False
$ws.Range("D64").Value = $True
$ws.Range("D65").Value = $True
$ws.Range("D66").Value = $False
$ws.Range("D69").Value = "[6 42 127; 41 76 247; 102 59 231; 162 41 216; 222 24 200; 255 192 203] / 255"
$ws.Range("D70").Value = 50
$ws.Range("D71").Value = 50
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 500
$ws.Range("D75").Value = 800
$ws.Range("D76").Value = 8
$ws.Range("D77").Value = 8
$ws.Range("D78").Value = 0.008
$ws.Range("D79").Value = 2000
$ws.Range("D80").Value = 4000
$ws.Range("D83").Value = "agentControl_Adam"
$ws.Range("D84").Value = "findNeighborhood_fixedRadius"

$ws.Range("G2").Select() | Out-Null
